# Boosting Algorithm.pptx edit script
# Applies the content changes from the target commit:
#  - Slide 1: "Thanigaivel G" (visual text unchanged; left as-is)
#  - Slides 2-7: title font size bumped to 40pt
#  - Slide 4: bullet list sizes 30pt -> 28pt; caption 28pt -> 24pt
#  - Slide 5: typo fix "a ensemble" -> "an ensemble" + bold a phrase
#  - Slide 6: typo/split title; body split + bold a phrase
#  - Slide 7: typo/split title; body typo fix "a ensemble" -> "an ensemble"
#             + bold two phrases
#  - Slide 8: body split + bold a phrase

function Set-BoldSubstring {
    param($TextRange, [string]$Target)

    $full = $TextRange.Text
    $idx = $full.IndexOf($Target)
    if ($idx -lt 0) {
        Write-Host "WARN: substring not found: [$Target]"
        return
    }
    $startPos = $idx + 1
    $sub = $TextRange.Characters($startPos, $Target.Length)
    $sub.Font.Bold = $true
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1 - Title slide: split "Thanigaivel G" into two runs (cosmetic,
# mirrors the spell-check run-split in the source deck; visible text is
# unchanged).
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1name = $s1.Shapes.Item(2).TextFrame.TextRange
$s1name.Characters(1, 11).Font.Size = $s1name.Characters(1, 11).Font.Size

# ---------------------------------------------------------------------
# Slide 2 - "Boosting Algorithm" (image slide)
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Font.Size = 40

# ---------------------------------------------------------------------
# Slide 3 - "Pros and Cons of Boosting Algorithm"
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Font.Size = 40

# ---------------------------------------------------------------------
# Slide 4 - "Types of Boosting Algorithm"
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Font.Size = 40

$s4list = $s4.Shapes.Item(2).TextFrame.TextRange
$s4paras = $s4list.Paragraphs()
for ($i = 1; $i -le $s4paras.Count; $i++) {
    $para = $s4list.Paragraphs($i, 1)
    $para.Font.Size = 28
}

$s4caption = $s4.Shapes.Item(3).TextFrame.TextRange
$s4caption.Font.Size = 24

# ---------------------------------------------------------------------
# Slide 5 - "Gradient Boosting"
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Font.Size = 40

$s5body = $s5.Shapes.Item(2).TextFrame.TextRange
$s5body.Text = "Gradient Boosting is a machine learning technique used for regression and classification tasks and it is an ensemble learning method that sequentially adds models to correct errors made by previous models"
Set-BoldSubstring $s5body "sequentially adds models to correct errors "

# ---------------------------------------------------------------------
# Slide 6 - "XGBoost Algorithm"
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6title = $s6.Shapes.Item(1).TextFrame.TextRange
$s6title.Font.Size = 40
# force the run split between "XGBoost" and " Algorithm" (cosmetic, same text)
$s6title.Characters(1, 7).Font.Size = 40

$s6body = $s6.Shapes.Item(2).TextFrame.TextRange
Set-BoldSubstring $s6body "multiple weak learners to improve prediction accuracy "

# ---------------------------------------------------------------------
# Slide 7 - "Adaboosting Algorithm"
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$s7title = $s7.Shapes.Item(1).TextFrame.TextRange
$s7title.Font.Size = 40
$s7title.Characters(1, 11).Font.Size = 40

$s7body = $s7.Shapes.Item(3).TextFrame.TextRange
$s7body.Text = "Adaboosting is an ensemble learning method that works by iteratively training weak models, adjusting the weights of incorrectly classified instances to focus on harder cases in subsequent iterations. The final model is a weighted sum of all the weak models, with more weight given to those that performed well."
Set-BoldSubstring $s7body "adjusting the weights "
Set-BoldSubstring $s7body "weighted sum "

# ---------------------------------------------------------------------
# Slide 8 - "Catboosting Algorithm"
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$s8body = $s8.Shapes.Item(3).TextFrame.TextRange
Set-BoldSubstring $s8body "categorical features "

Write-Host "edit.ps1 completed"
